# Update cryptocurrency price/volume data to reflect latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'330.72"
$ws.Range("E2").Value = "'-0.23%"
$ws.Range("D3").Value = "'41.61"
$ws.Range("E3").Value = "'0.92%"
$ws.Range("D4").Value = "'5.698"
$ws.Range("E4").Value = "'0.21%"
$ws.Range("D5").Value = "'0.08392"
$ws.Range("E5").Value = "'3.93%"
$ws.Range("D6").Value = "'8.820"
$ws.Range("E6").Value = "'0.80%"
$ws.Range("D7").Value = "'2.008"
$ws.Range("E7").Value = "'-1.36%"
$ws.Range("D8").Value = "'4.476"
$ws.Range("E8").Value = "'-1.34%"
$ws.Range("D10").Value = "'0.9238"
$ws.Range("E10").Value = "'-0.10%"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("E11").Value = "'1.27%"
$ws.Range("D12").Value = "'0.1973"
$ws.Range("E12").Value = "'1.47%"
$ws.Range("D13").Value = "'0.09542"
$ws.Range("E13").Value = "'1.58%"
$ws.Range("D14").Value = "'0.03884"
$ws.Range("E14").Value = "'3.47%"
$ws.Range("E15").Value = "'0.69%"
$ws.Range("D16").Value = "'0.001298"
$ws.Range("E16").Value = "'-0.10%"
$ws.Range("E17").Value = "'-3.23%"
$ws.Range("D18").Value = "'3.424"
$ws.Range("E18").Value = "'1.81%"
$ws.Range("E19").Value = "'0.67%"
$ws.Range("D20").Value = "'8.869"
$ws.Range("E20").Value = "'1.55%"
$ws.Range("E21").Value = "'-3.88%"
$ws.Range("D23").Value = "'0.04413"
$ws.Range("E23").Value = "'-0.62%"
$ws.Range("D24").Value = "'0.001273"
$ws.Range("E24").Value = "'0.91%"
$ws.Range("D25").Value = "'0.004401"
$ws.Range("E25").Value = "'2.49%"
$ws.Range("E26").Value = "'-4.20%"
$ws.Range("D27").Value = "'0.0003992"
$ws.Range("E27").Value = "'-0.03%"
$ws.Range("D39").Value = "'0.02863"
$ws.Range("E39").Value = "'-0.27%"
$ws.Range("E40").Value = "'0.64%"
$ws.Range("E41").Value = "'2.11%"
$ws.Range("D42").Value = "'0.1433"
$ws.Range("E42").Value = "'0.90%"
$ws.Range("D43").Value = "'0.008978"
$ws.Range("E43").Value = "'-9.75%"
$ws.Range("E44").Value = "'-2.53%"
$ws.Range("D45").Value = "'0.01167"
$ws.Range("E45").Value = "'-1.26%"
$ws.Range("D46").Value = "'0.00006938"
$ws.Range("E46").Value = "'2.30%"
$ws.Range("E47").Value = "'-0.18%"
$ws.Range("D48").Value = "'0.003468"
$ws.Range("E48").Value = "'14.69%"
$ws.Range("E49").Value = "'-0.24%"
$ws.Range("E50").Value = "'-0.18%"
$ws.Range("E51").Value = "'-0.18%"
